$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.020126409828662872
$ws.Range("C2").Value = 0.00830205250531435
$ws.Range("D2").Value = 0.0065217469818890095
$ws.Range("E2").Value = 0.003844057209789753
$ws.Range("F2").Value = 0.0000051219753913755994
$ws.Range("J2").Value = 0.12627831101417542
$ws.Range("K2").Value = 1.425947666168213
